$d = $word.ActiveDocument

function New-RunXml($text, $preserve) {
    $escaped = $text.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
    $spaceAttr = ""
    if ($preserve) { $spaceAttr = ' xml:space="preserve"' }
    return '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="da-DK"/></w:rPr><w:t' + $spaceAttr + '>' + $escaped + '</w:t></w:r>'
}

function New-PackageXml($runsXml) {
    return @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
'@ + $runsXml + @'
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@
}

function Replace-ParagraphContent($paragraph, $segments) {
    $pStart = $paragraph.Range.Start
    $oldEnd = $paragraph.Range.End

    $runsXml = ""
    $insertedLen = 0
    foreach ($seg in $segments) {
        $text = $seg[0]
        $preserve = $seg[1]
        $runsXml += New-RunXml $text $preserve
        $insertedLen += $text.Length
    }

    $xml = New-PackageXml $runsXml
    $insertionPoint = $d.Range($pStart, $pStart)
    $insertionPoint.InsertXML($xml)

    $newOldStart = $pStart + $insertedLen
    $newOldEnd = $oldEnd + $insertedLen
    $oldRange = $d.Range($newOldStart, $newOldEnd)
    $oldRange.Delete()
}

# --- Paragraph 2 (first content paragraph) -------------------------------
$para2 = $d.Paragraphs.Item(2)
$para2Segments = @(
    , @("Borgerne", $false)
    , @(" i Region Nordjylland klager over den lange ventetid på telefonlinjen, når de bestiller Flexsygehustaxa. ", $true)
    , @("Derfor vælger ", $true)
    , @("borgerne", $false)
    , @(" i stedet at ringe til lægerne, så deres ", $true)
    , @("lægesekretærer", $false)
    , @(" kan bestille taxaerne for dem. Dette resulterer dermed", $true)
    , @(" i", $true)
    , @(",", $false)
    , @(" at ventetiden pr. bestilling er kortere", $true)
    , @(", da ", $true)
    , @("lægesekretærer", $false)
    , @("ne kan ringe via et internt nummer.", $false)
    , @(" ", $true)
    , @("D", $false)
    , @("og ", $true)
    , @("tilføjer", $false)
    , @(" det også ", $true)
    , @("en ny opgave hos ", $true)
    , @("lægesekretærer", $false)
    , @("n", $false)
    , @("e, som kunne undgås hvis ", $true)
    , @("borgeren", $false)
    , @(" havde bestilt taxaen selv", $true)
    , @(".", $false)
)
Replace-ParagraphContent $para2 $para2Segments

# --- Paragraph 3 (second content paragraph) -------------------------------
$para3 = $d.Paragraphs.Item(3)
$para3Segments = @(
    , @("Hvordan kan man få ", $true)
    , @("borgerne", $false)
    , @(" til selv at bestille deres ", $true)
    , @("F", $false)
    , @("lexsygehustaxa uden hjælp fra ", $true)
    , @("lægesekretærer", $false)
    , @("ne og uden den unødvendige lange ventetid", $false)
    , @(", og kan man ellers", $false)
    , @(" formindske lægesekretærernes ekstra arbejde?", $true)
)
Replace-ParagraphContent $para3 $para3Segments

Write-Host "Done. Para2: $($para2.Range.Text)"
Write-Host "Done. Para3: $($para3.Range.Text)"
